$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.065365076065063
$ws.Range("B1").Value = 1.217955946922302
$ws.Range("C1").Value = 1.530105948448181
$ws.Range("D1").Value = 3.132083892822266
$ws.Range("E1").Value = 4.253100872039795
